$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add O4 = 2021, matching the format of N4 (header year row)
$ws.Range("O4").Value = 2021
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)

# Add O5 = 1.5020015556876996, matching the format of N5 (data row)
$ws.Range("O5").Value = 1.5020015556876996
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Move the active selection to Q5 (matches recorded diff)
$ws.Range("Q5").Select()
